$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -12.214
$ws.Range("C9").Value = -10.3342
$ws.Range("A11").Value = -21.75069999999999
$ws.Range("B11").Value = 5.545100000000002
$ws.Range("A12").Value = -21.36579999999999
$ws.Range("C13").Value = -12.9302
$ws.Range("C14").Value = -13.90509999999999
$ws.Range("A15").Value = -21.61610000000001
$ws.Range("C19").Value = -12.50490000000002
$ws.Range("C21").Value = -12.5892
$ws.Range("C22").Value = -11.4875
$ws.Range("B23").Value = 8.818299999999997
$ws.Range("C24").Value = -12.72249999999999
$ws.Range("C26").Value = -11.90829999999999
$ws.Range("A27").Value = -21.8916
$ws.Range("A28").Value = -21.80420000000001
$ws.Range("B28").Value = 5.881699999999999
$ws.Range("A31").Value = -21.5195
$ws.Range("A32").Value = -21.11929999999997
$ws.Range("B32").Value = 5.854399999999997
$ws.Range("B34").Value = 9.469200000000003
$ws.Range("A36").Value = -21.27589999999998
$ws.Range("B36").Value = 5.165999999999998
$ws.Range("B37").Value = 9.033399999999999
$ws.Range("A38").Value = -20.04109999999998
$ws.Range("C38").Value = -11.9182
$ws.Range("C41").Value = -12.72930000000001
$ws.Range("B42").Value = 10.3179
$ws.Range("A46").Value = -21.82080000000001
$ws.Range("B49").Value = 4.9772
$ws.Range("C52").Value = -10.9439
$ws.Range("A54").Value = -21.77040000000001
$ws.Range("B54").Value = 4.288900000000001
$ws.Range("A55").Value = -21.83620000000001
$ws.Range("A56").Value = -21.96339999999999
$ws.Range("C56").Value = -12.5282
$ws.Range("A67").Value = -21.60569999999997
$ws.Range("A69").Value = -21.77169999999998
$ws.Range("C71").Value = -12.7926
$ws.Range("A72").Value = -21.93209999999999
$ws.Range("C72").Value = -12.41379999999999
$ws.Range("A73").Value = -19.35470000000002
$ws.Range("B78").Value = 9.919900000000005
$ws.Range("C78").Value = -11.59550000000001
$ws.Range("B80").Value = 9.139799999999997
$ws.Range("A83").Value = -21.76099999999999
$ws.Range("C83").Value = -12.2299
$ws.Range("C85").Value = -13.64569999999999
$ws.Range("A86").Value = -21.8506
$ws.Range("C86").Value = -12.8688
$ws.Range("C90").Value = -10.0952
$ws.Range("A91").Value = -20.7256
$ws.Range("A93").Value = -21.4098
$ws.Range("C96").Value = -10.2097
$ws.Range("B97").Value = 6.337899999999996
$ws.Range("A99").Value = -21.85840000000001
$ws.Range("B99").Value = 5.057699999999998
$ws.Range("B100").Value = 4.482400000000002
$ws.Range("B101").Value = 4.676499999999999
$ws.Range("C103").Value = -13.30099999999999
$ws.Range("A104").Value = -21.54969999999999
$ws.Range("A105").Value = -19.94969999999999
